$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.602.57'
$ws.Range("E2").Value = '  +4.35%  '
$ws.Range("D3").Value = '3.058.52'
$ws.Range("E3").Value = '  +3.21%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.996'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '551.77'
$ws.Range("E5").Value = '  +1.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.66'
$ws.Range("E6").Value = '  +10.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("D8").Value = '3.058.51'
$ws.Range("E8").Value = '  +3.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.495'
$ws.Range("E9").Value = '  +1.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.95'
$ws.Range("E10").Value = '  +19.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.152'
$ws.Range("E11").Value = '  +5.30%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.453'
$ws.Range("E12").Value = '  +3.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000225'
$ws.Range("E13").Value = '  +3.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.83'
$ws.Range("E14").Value = '  +4.04%  '
$ws.Range("D15").Value = '3.506.57'
$ws.Range("E15").Value = '  +1.81%  '
$ws.Range("D16").Value = '63.239.43'
$ws.Range("E16").Value = '  +3.78%  '
$ws.Range("D17").Value = '3.035.80'
$ws.Range("E17").Value = '  +2.48%  '
$ws.Range("E18").Value = '  -1.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.70'
$ws.Range("E19").Value = '  +2.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '477.66'
$ws.Range("E20").Value = '  +3.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.55'
$ws.Range("E21").Value = '  +4.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.666'
$ws.Range("E22").Value = '  +1.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.40'
$ws.Range("E23").Value = '  +7.68%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.00'
$ws.Range("E24").Value = '  +10.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.30'
$ws.Range("E25").Value = '  +1.46%  '
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.76'
$ws.Range("E27").Value = '  +3.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.85'
$ws.Range("E28").Value = '  +4.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.05'
$ws.Range("E29").Value = '  +9.71%  '
$ws.Range("E30").Value = '  -0.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.99'
$ws.Range("E31").Value = '  +3.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.14'
$ws.Range("E32").Value = '  +2.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.42'
$ws.Range("E33").Value = '  +7.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.69'
$ws.Range("E34").Value = '  +6.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '55.11'
$ws.Range("E35").Value = '  +1.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.00'
$ws.Range("E36").Value = '  +3.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '458.89'
$ws.Range("E37").Value = '  +3.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0824'
$ws.Range("E38").Value = '  +5.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0400'
$ws.Range("E39").Value = '  +7.37%  '
$ws.Range("D40").Value = '2.986.04'
$ws.Range("E40").Value = '  -4.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.116'
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.24'
$ws.Range("E42").Value = '  +3.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.67'
$ws.Range("E43").Value = '  +14.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '27.56'
$ws.Range("E44").Value = '  +8.98%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.253'
$ws.Range("E45").Value = '  +6.81%  '
$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.04'
$ws.Range("E47").Value = '  +6.47%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.111'
$ws.Range("E48").Value = '  +3.99%  '
$ws.Range("D49").Value = '0.0₃0514'
$ws.Range("E49").Value = '  +8.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '115.77'
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.05'
$ws.Range("E51").Value = '  +4.35%  '
